$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8108.478
$ws.Range("I40").Value = 4399.6
$ws.Range("J40").Value = 9138.723
$ws.Range("K40").Value = 4399.6
$ws.Range("L40").Value = 9138.723
$ws.Range("M40").Value = -4224.6
$ws.Range("N40").Value = -9488.723

$ws.Range("H53").Value = 284.26666
$ws.Range("I53").Value = 315.2
$ws.Range("K53").Value = 315.2
$ws.Range("M53").Value = 321.8

$ws.Range("H69").Value = 8011.7144
$ws.Range("I69").Value = 8007.6665
$ws.Range("J69").Value = 8014.75
$ws.Range("K69").Value = 24022.9995
$ws.Range("L69").Value = 24044.25
$ws.Range("M69").Value = -23148.9995
$ws.Range("N69").Value = -25792.25

$ws.Range("H72").Value = 8011.7144
$ws.Range("I72").Value = 8007.6665
$ws.Range("J72").Value = 8014.75
$ws.Range("K72").Value = 72068.9985
$ws.Range("L72").Value = 72132.75
$ws.Range("M72").Value = -67700.9985
$ws.Range("N72").Value = -80868.75

$ws.Range("H74").Value = 5375
$ws.Range("I74").Value = 5375
$ws.Range("K74").Value = 5375
$ws.Range("M74").Value = -4439

$ws.Range("H77").Value = 5375
$ws.Range("I77").Value = 5375
$ws.Range("K77").Value = 26875
$ws.Range("M77").Value = -22195

$ws.Range("H100").Value = 2501
$ws.Range("I100").Value = 2501
$ws.Range("K100").Value = 2501
$ws.Range("M100").Value = -1960

$ws.Range("H118").Value = 279.5
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H129").Value = 3075.2
$ws.Range("J129").Value = 4999.5
$ws.Range("L129").Value = 14998.5
$ws.Range("N129").Value = -24998.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 521.5714
$ws.Range("I4").Value = 488.2
$ws.Range("K4").Value = 488.2
$ws.Range("M4").Value = -372.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 388.22223
$ws.Range("I22").Value = 388.22223
$ws.Range("K22").Value = 388.22223
$ws.Range("M22").Value = -215.22223

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 898.1667
$ws.Range("I7").Value = 898.1667
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 898.1667
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -785.1667
$ws.Range("N7").ClearContents()

$ws.Range("H12").Value = 394
$ws.Range("J12").Value = 350
$ws.Range("L12").Value = 350
$ws.Range("N12").Value = -690

$ws.Range("H13").Value = 1850
$ws.Range("I13").Value = 1325
$ws.Range("K13").Value = 1325
$ws.Range("M13").Value = -1186

$ws.Range("H19").Value = 96
$ws.Range("I19").Value = 93.8
$ws.Range("J19").Value = 98.75
$ws.Range("K19").Value = 93.8
$ws.Range("L19").Value = 98.75
$ws.Range("M19").Value = 76.2
$ws.Range("N19").Value = -438.75

$ws.Range("H23").Value = 883.3333
$ws.Range("J23").Value = 850
$ws.Range("L23").Value = 850
$ws.Range("N23").Value = -1330

$ws.Range("H24").Value = 96
$ws.Range("I24").Value = 93.8
$ws.Range("J24").Value = 98.75
$ws.Range("K24").Value = 93.8
$ws.Range("L24").Value = 98.75
$ws.Range("M24").Value = 76.2
$ws.Range("N24").Value = -438.75

$ws.Range("H27").Value = 883.3333
$ws.Range("J27").Value = 850
$ws.Range("L27").Value = 850
$ws.Range("N27").Value = -1234

$ws.Range("H31").Value = 2979.3333
$ws.Range("I31").Value = 1486.2858
$ws.Range("J31").Value = 4285.75
$ws.Range("K31").Value = 1486.2858
$ws.Range("L31").Value = 4285.75
$ws.Range("M31").Value = -1191.2858
$ws.Range("N31").Value = -4875.75

$ws.Range("H34").Value = 2979.3333
$ws.Range("I34").Value = 1486.2858
$ws.Range("J34").Value = 4285.75
$ws.Range("K34").Value = 1486.2858
$ws.Range("L34").Value = 4285.75
$ws.Range("M34").Value = -1284.2858
$ws.Range("N34").Value = -4689.75

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1063.4286
$ws.Range("I23").Value = 886
$ws.Range("J23").Value = 1300
$ws.Range("K23").Value = 2658
$ws.Range("L23").Value = 3900
$ws.Range("M23").Value = -2423
$ws.Range("N23").Value = -4370

$ws.Range("H24").Value = 5000
$ws.Range("I24").Value = 5000
$ws.Range("K24").Value = 15000
$ws.Range("M24").Value = -14770

$ws.Range("H103").Value = 1781.25
$ws.Range("I103").Value = 998.5
$ws.Range("J103").Value = 2172.625
$ws.Range("K103").Value = 2995.5
$ws.Range("L103").Value = 6517.875
$ws.Range("M103").Value = -2116.5
$ws.Range("N103").Value = -8275.875

$ws.Range("H113").Value = 684.7143
$ws.Range("J113").Value = 758
$ws.Range("L113").Value = 2274
$ws.Range("N113").Value = -6614

$ws.Range("H133").Value = 9749.5
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 600
$ws.Range("I134").Value = 600
$ws.Range("K134").Value = 1800
$ws.Range("M134").Value = 3270

$ws.Range("H136").Value = 1278
$ws.Range("I136").Value = 1278
$ws.Range("K136").Value = 3834
$ws.Range("M136").Value = 1266

$ws.Range("H138").Value = 2471
$ws.Range("I138").Value = 2404.875
$ws.Range("K138").Value = 7214.625
$ws.Range("M138").Value = -2074.625

$ws.Range("H139").Value = 2199.8
$ws.Range("I139").Value = 1874.75
$ws.Range("K139").Value = 5624.25
$ws.Range("M139").Value = -484.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 250916.5
$ws.Range("I3").Value = 999.75
$ws.Range("J3").Value = 750750
$ws.Range("K3").Value = 999.75
$ws.Range("L3").Value = 750750
$ws.Range("M3").Value = -883.75
$ws.Range("N3").Value = -750982

$ws.Range("H10").Value = 5003
$ws.Range("I10").Value = 5003
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 5003
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -4834
$ws.Range("N10").ClearContents()

$ws.Range("H20").Value = 13121
$ws.Range("I20").Value = 8005
$ws.Range("J20").Value = 14400
$ws.Range("K20").Value = 8005
$ws.Range("L20").Value = 14400
$ws.Range("M20").Value = -7760
$ws.Range("N20").Value = -14890

$ws.Range("H24").Value = 13000
$ws.Range("J24").Value = 13000
$ws.Range("L24").Value = 13000
$ws.Range("N24").Value = -13346

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4746
$ws.Range("I61").Value = 4746
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4746
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4544
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 4746
$ws.Range("I113").Value = 4746
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4746
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2576
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2887

$ws.Range("H9").Value = 3006
$ws.Range("I9").Value = 3006
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 3006
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("M9").Value = -2866

$ws.Range("H107").Value = 1716.0834
$ws.Range("I107").Value = 859.3
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 2577.9
$ws.Range("L107").Value = 18000
$ws.Range("M107").Value = -657.8999999999996
$ws.Range("N107").Value = -21840
